# Update "想去人数" (F column) values across sheets to match newly scraped data.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1188
$ws1.Range("F4").Value = 0
$ws1.Range("F5").Value = 5007
$ws1.Range("F8").Value = 239
$ws1.Range("F11").Value = 0
$ws1.Range("F12").Value = 0

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 0
$ws2.Range("F4").Value = 0
$ws2.Range("F5").Value = 6
$ws2.Range("F6").Value = 0

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 349
$ws4.Range("F5").Value = 14
$ws4.Range("F7").Value = 5007
$ws4.Range("F8").Value = 523
$ws4.Range("F10").Value = 0
$ws4.Range("F12").Value = 527
$ws4.Range("F13").Value = 86
$ws4.Range("F14").Value = 0
$ws4.Range("F15").Value = 3
$ws4.Range("F16").Value = 664
